$wb = $excel.ActiveWorkbook

# --- Rename and reorder sheets -------------------------------------------
# Old Sheet1 (admin success data) becomes "LoginSuccess"
# Old Sheet2 (admin fail / PM data) becomes "LoginFail"
# A brand-new sheet "LoginData" is inserted in front of them and becomes active.
#
# NOTE: worksheet handles returned by this host resolve by sheet *position*,
# not identity - once a sheet is inserted/removed/reordered, any previously
# captured handle can silently start pointing at a different sheet. So every
# handle below is re-fetched by NAME immediately before it is used.
$wb.Worksheets.Item(1).Name = "LoginSuccess"
$wb.Worksheets.Item(2).Name = "LoginFail"

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("LoginSuccess"))
$newSheet.Name = "LoginData"

# --- Populate the new LoginData sheet -------------------------------------
$wsData = $wb.Worksheets.Item("LoginData")

$wsData.Range("A1").Value2 = "email"
$wsData.Range("B1").Value2 = "password"

# Column by column so new shared-string entries are appended in the same
# order the original author typed them in (email column first, then
# password column).
$wsData.Range("A2").Value2 = "admin@example.com"
$wsData.Range("A3").Value2 = "admin123@example.com"
$wsData.Range("A4").Value2 = "employee123@example.com"

$wsData.Range("B2").Value2 = "'123456"
$wsData.Range("B3").Value2 = "'1234567"
$wsData.Range("B4").Value2 = "'1234568"

# Hyperlinks for the e-mail column
$wsData.Hyperlinks.Add($wsData.Range("A2"), "mailto:admin@example.com") | Out-Null
$wsData.Range("A2").Style = "Hyperlink"
$wsData.Hyperlinks.Add($wsData.Range("A3"), "mailto:admin123@example.com") | Out-Null
$wsData.Range("A3").Style = "Hyperlink"
$wsData.Hyperlinks.Add($wsData.Range("A4"), "mailto:employee123@example.com") | Out-Null
$wsData.Range("A4").Style = "Hyperlink"

# Column widths matching the other sheets
$wsData.Columns.Item(1).ColumnWidth = 26.592447916666668
$wsData.Columns.Item(2).ColumnWidth = 12.166666666666666

# --- LoginFail sheet: normalise the highlighted row's alignment ----------
$wsFail = $wb.Worksheets.Item("LoginFail")
$wsFail.Range("A2:C2").HorizontalAlignment = -4131  # xlLeft
$wsFail.Range("A2:C2").VerticalAlignment = -4160    # xlTop

# --- Selections / active-tab state ----------------------------------------
# Each sheet's selection is only recorded for whichever sheet is active at
# the moment .Select() runs, so activate-then-select one at a time, doing
# the sheet that should stay "on top" (LoginData) last.
$wb.Worksheets.Item("LoginFail").Activate()
$wb.Worksheets.Item("LoginFail").Range("B35").Select() | Out-Null

$wb.Worksheets.Item("LoginSuccess").Activate()
$wb.Worksheets.Item("LoginSuccess").Range("A4:XFD4").Select() | Out-Null

$wb.Worksheets.Item("LoginData").Activate()
$wb.Worksheets.Item("LoginData").Range("C1:F1048576").Select() | Out-Null
